# Applies the "Add link to final project" edit:
#   - Adds two new bullet points under "Next steps:" after the existing
#     "Kleinigkeiten ueberarbeiten (...)" bullet:
#       * "Conclusio in Form eines Sheets einfuegen -> vielleicht mit
#         Anmerkungen oder auch noch einfachen Stichpunkten"
#       * "Correlation explorer in die Story einbinden "
#   - The "_GoBack" bookmark, which previously sat at the end of the
#     "Kleinigkeiten ueberarbeiten" bullet, moves to sit between
#     "...Anmerkungen" and " oder auch noch einfachen Stichpunkten" in the
#     new first bullet.

$d = $word.ActiveDocument

# Locate the "Kleinigkeiten ueberarbeiten ..." bullet (the last bullet under
# "Next steps:") by scanning paragraphs for a distinctive snippet of its
# text, so the script does not depend on a hard-coded paragraph index.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Kleinigkeiten*") {
        $anchorIndex = $i
        break
    }
}
$anchorPara = $d.Paragraphs.Item($anchorIndex)

# The existing "_GoBack" bookmark sits at the end of that bullet; remove it
# here so it can be re-created at its new location further down.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert a new list paragraph right after the "Kleinigkeiten..." bullet; it
# inherits the list/paragraph style (Listenabsatz, numId 6) automatically.
$anchorPara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($anchorIndex + 1)
$newPara1.Range.Text = "Conclusio in Form eines Sheets einfügen -> vielleicht mit Anmerkungen oder auch noch einfachen Stichpunkten"

# Re-create the "_GoBack" bookmark between "...Anmerkungen" and " oder auch
# noch einfachen Stichpunkten" in the new paragraph.
$bmRange = $newPara1.Range.Duplicate
$bmRange.Find.Execute("Anmerkungen") | Out-Null
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Insert a second new list paragraph after the first new one.
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($anchorIndex + 2)
$newPara2.Range.Text = "Correlation explorer in die Story einbinden "
